$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Remove the first worker (ARIEL DE JESUS YEPES CASTRO / 1047403115,
#    period 2104) which occupied row 16. Deleting the whole row shifts
#    every following row (table rows + the signature block) up by one,
#    which is exactly what the target workbook shows (table now ends at
#    row 31 instead of 32, signature rows become 36/37 instead of 37/38).
# ------------------------------------------------------------------
$ws.Range("B16").EntireRow.Delete()

# ------------------------------------------------------------------
# 2) Refresh the account-statement data table (rows 16-31): two workers
#    (EVER BARRIOS PUELLO / LINA MARGARITA RANGEL MONTES), periods 2201
#    through 2208 in ascending order, alternating between the two
#    workers, with updated "Valor Mora" (F) / "Salario Basico" (G)
#    amounts.
# ------------------------------------------------------------------
$doc1 = "70526731"
$name1 = "EVER BARRIOS PUELLO"
$doc2 = "1143402322"
$name2 = "LINA MARGARITA RANGEL MONTES"

$rows = @(
  @{ Row=16; Doc=$doc1; Name=$name1; Period="2201"; F=38000; G=1472625 },
  @{ Row=17; Doc=$doc2; Name=$name2; Period="2201"; F=53480; G=1337000 },
  @{ Row=18; Doc=$doc1; Name=$name1; Period="2202"; F=38000; G=1472625 },
  @{ Row=19; Doc=$doc2; Name=$name2; Period="2202"; F=53480; G=1337000 },
  @{ Row=20; Doc=$doc1; Name=$name1; Period="2203"; F=67320; G=1472625 },
  @{ Row=21; Doc=$doc2; Name=$name2; Period="2203"; F=53480; G=1337000 },
  @{ Row=22; Doc=$doc1; Name=$name1; Period="2204"; F=67320; G=1472625 },
  @{ Row=23; Doc=$doc2; Name=$name2; Period="2204"; F=53480; G=1337000 },
  @{ Row=24; Doc=$doc1; Name=$name1; Period="2205"; F=67320; G=1472625 },
  @{ Row=25; Doc=$doc2; Name=$name2; Period="2205"; F=53480; G=1337000 },
  @{ Row=26; Doc=$doc1; Name=$name1; Period="2206"; F=67320; G=1472625 },
  @{ Row=27; Doc=$doc2; Name=$name2; Period="2206"; F=53480; G=1337000 },
  @{ Row=28; Doc=$doc1; Name=$name1; Period="2207"; F=67320; G=1472625 },
  @{ Row=29; Doc=$doc2; Name=$name2; Period="2207"; F=53480; G=1337000 },
  @{ Row=30; Doc=$doc1; Name=$name1; Period="2208"; F=43197; G=1472625 },
  @{ Row=31; Doc=$doc2; Name=$name2; Period="2208"; F=39219; G=1337000 }
)

foreach ($r in $rows) {
  $n = $r.Row
  $ws.Range("B$n").Value = "CC"
  $ws.Range("C$n").Value = $r.Doc
  $ws.Range("D$n").Value = $r.Name
  $ws.Range("E$n").Value = $r.Period
  $ws.Range("F$n").Value = $r.F
  $ws.Range("G$n").Value = $r.G
}

# ------------------------------------------------------------------
# 3) Refresh the summary fields above the table: total "Valor Mora"
#    (E11), worker count (C13) and period count (F13).
# ------------------------------------------------------------------
$ws.Range("E11").Value = 869376
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 8
